$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 14: Hyōgo -> Hyogo ; Visit_hx No -> Unknown
$ws.Range("A14").Value = "Hyogo"
$ws.Range("G14").Value = "Unknown"

# Row 21: Kōchi -> Kochi ; Visit_hx No -> Unknown
$ws.Range("A21").Value = "Kochi"
$ws.Range("G21").Value = "Unknown"

# Row 23: Kyōto -> Kyoto ; Visit_hx Yes -> Unknown
$ws.Range("A23").Value = "Kyoto"
$ws.Range("G23").Value = "Unknown"

# Row 31: Ōita -> Oita ; Visit_hx No -> Unknown
$ws.Range("A31").Value = "Oita"
$ws.Range("G31").Value = "Unknown"

# Row 34: Ōsaka -> Osaka ; Visit_hx Yes -> Unknown
$ws.Range("A34").Value = "Osaka"
$ws.Range("G34").Value = "Unknown"

# Row 42: Tōkyō -> Tokyo ; Visit_hx Yes -> Unknown
$ws.Range("A42").Value = "Tokyo"
$ws.Range("G42").Value = "Unknown"
